$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores crypto prices as plain text (e.g. "59.181.17", a
# thousands-dotted value, or "0.0₃0778", a subscript-notation value) -
# never as real numbers. Assigning a numeric-looking string straight to
# .Value lets Excel auto-coerce it into a Number, which both changes the
# cell type and can introduce floating-point rounding noise. Flip the
# column to Text first so the new values land as text, then restore the
# "Normal" cell style so no stray number-format/style is left behind -
# matching source data, which never touches cell formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '59.181.17'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '2.525.44'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '537.66'
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('D6').Value = '138.47'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '0.566'
$ws.Range('E8').Value = '  +0.34%  '
$ws.Range('D9').Value = '2.523.29'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('D11').Value = '0.159'
$ws.Range('D12').Value = '5.37'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('D13').Value = '0.347'
$ws.Range('E13').Value = '  -2.34%  '
$ws.Range('D14').Value = '2.961.68'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').Value = '23.22'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '59.017.53'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '2.527.19'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').Value = '11.10'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').Value = '4.30'
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('D21').Value = '325.80'
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = '5.90'
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('D24').Value = '65.69'
$ws.Range('E24').Value = '  +5.42%  '
$ws.Range('D25').Value = '0.425'
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').Value = '7.66'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '6.76'
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0778'
$ws.Range('E30').Value = '  +1.21%  '
$ws.Range('D31').Value = '1.78'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '169.69'
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('E33').Value = '  +7.36%  '
$ws.Range('E34').Value = '  +2.84%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = '18.57'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('D37').Value = '4.12'
$ws.Range('E37').Value = '  -2.07%  '
$ws.Range('D38').Value = '1.57'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '36.74'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').Value = '0.828'
$ws.Range('E40').Value = '  +3.35%  '
$ws.Range('D41').Value = '3.64'
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('D42').Value = '284.75'
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('D43').Value = '5.26'
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = '130.88'
$ws.Range('E45').Value = '  +7.76%  '
$ws.Range('E46').Value = '  +1.61%  '
$ws.Range('D47').Value = '10.86'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').Value = '0.0933'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('E51').Value = '  +0.19%  '

$ws.Range("D2:D51").Style = "Normal"
